# Added 200 Packet run MW 1D
# Populates the "200 Packets" (Acc/Loss/Time) results block in columns Z:AB
# of the "MW CNN 1D" sheet, mirroring the existing 298/290/270/250/230/210
# packet blocks, plus a stray label cell and the corresponding average rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MW CNN 1D")
$ws.Activate()

# --- Header rows -----------------------------------------------------
# Row 1: group label "200 Packets" over the new block (merged-looking header
# sits in the middle column, AA, like the other blocks sit over their
# middle column).
$ws.Range("AA1").Value = "200 Packets"

# Row 2: per-column sub headers (Acc / Loss / Time)
$ws.Range("Z2").Value = "Acc"
$ws.Range("AA2").Value = "Loss"
$ws.Range("AB2").Value = "Time"

# --- Data rows (3-51): Acc / Loss / Time for the 200 Packets trials ---
$rows = @(
    @(3, 88.647520542144704, 0.29172533327497502, 471.59354948997498),
    @(4, 89.394551515579195, 0.25142665083347998, 461.119449615478),
    @(5, 89.271229505538898, 0.26566754305444901, 461.11307930946299),
    @(6, 89.688616991043006, 0.25564228762697799, 461.56454730033801),
    @(7, 88.019067049026404, 0.29991020372850602, 461.12673068046502),
    @(8, 89.000874757766695, 0.28511941850529499, 465.46656298637299),
    @(9, 88.730525970458899, 0.27735378177696901, 460.39963817596401),
    @(10, 89.247518777847205, 0.25306278556266498, 464.23108601570101),
    @(11, 89.079135656356797, 0.24944038304830901, 462.913558006286),
    @(12, 85.476791858673096, 0.32642584738585301, 463.599218606948),
    @(13, 88.699692487716604, 0.27395483164948198, 462.34886145591702),
    @(14, 89.1194522380828, 0.27031771345587702, 463.07760882377602),
    @(15, 89.847511053085299, 0.25591995566252901, 460.390878915786),
    @(16, 88.301277160644503, 0.26133078872799798, 462.88062715530299),
    @(17, 87.457019090652395, 0.29450335650544102, 462.04888916015602),
    @(18, 88.934475183486896, 0.25481036988731598, 461.48591017722998),
    @(19, 87.867289781570406, 0.31816218429759602, 463.88524746894802),
    @(20, 89.017480611801105, 0.257885398734092, 461.82305645942603),
    @(21, 88.799297809600802, 0.27025974393110003, 463.57051944732598),
    @(22, 89.233285188674898, 0.25176019515860398, 462.14640021324101),
    @(23, 89.572411775588904, 0.245837544843688, 463.55166697502102),
    @(24, 89.010363817214895, 0.27472634441592397, 459.59633493423399),
    @(25, 89.247518777847205, 0.26139393250298198, 461.20588707923798),
    @(26, 88.718664646148596, 0.27351009722826503, 465.39883661270102),
    @(27, 89.114707708358694, 0.25417985736045901, 462.37801456451399),
    @(28, 88.694953918457003, 0.293079361334737, 463.96956706047001),
    @(29, 89.143168926238999, 0.24245275793190299, 459.76927876472399),
    @(30, 89.000874757766695, 0.25922572803289601, 460.57213950157097),
    @(31, 88.481515645980807, 0.29369116100866999, 462.12938785552899),
    @(32, 89.181113243102999, 0.25030409383191898, 462.83479976654002),
    @(33, 87.864917516708303, 0.298987714733779, 462.400925159454),
    @(34, 88.777953386306706, 0.27118690311007498, 464.17543506622297),
    @(35, 88.815897703170705, 0.27288601945761498, 462.15470814704798),
    @(36, 89.306801557540894, 0.24464437900555899, 460.146674156188),
    @(37, 89.169257879257202, 0.25969924925655202, 461.79513716697602),
    @(38, 87.978750467300401, 0.31268606102309399, 461.53372931480402),
    @(39, 89.323407411575303, 0.25088632003663403, 462.360209465026),
    @(40, 87.532907724380493, 0.338108114561606, 460.29124689102099),
    @(41, 88.996136188506995, 0.26961057123740201, 463.76413369178698),
    @(42, 89.529728889465304, 0.24368283545649699, 463.56668210029602),
    @(43, 88.839614391326904, 0.25797176401899802, 462.25898408889702),
    @(44, 89.188230037689195, 0.26381842122476601, 460.543656826019),
    @(45, 89.188230037689195, 0.25237014731330198, 463.21752071380598),
    @(46, 87.891006469726506, 0.32553704761340202, 460.541975975036),
    @(47, 88.460171222686697, 0.27916198640133699, 463.69074869155799),
    @(48, 88.858586549758897, 0.27800204176438498, 463.50334715843201),
    @(49, 88.455426692962604, 0.27894910599501899, 463.32916593551602),
    @(50, 88.000094890594397, 0.28667592244733697, 462.64667034149102),
    @(51, 89.043563604354802, 0.24802463422834001, 462.66078042983997)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 26).Value = $r[1]   # Z  = Acc
    $ws.Cells.Item($rowNum, 27).Value = $r[2]   # AA = Loss
    $ws.Cells.Item($rowNum, 28).Value = $r[3]   # AB = Time
}

# Stray "v c" label that appears mid-table in the source workbook (row 40,
# column M) alongside the 270-Packets block.
$ws.Range("M40").Value = "v c"

# --- Average rows ------------------------------------------------------
$ws.Range("Z53").Formula = "=AVERAGE(Z3:Z51)"
$ws.Range("AA53").Formula = "=AVERAGE(AA3:AA51)"
$ws.Range("AB53").Formula = "=AVERAGE(AB3:AB51)"

$ws.Range("Z54").Formula = "=AVERAGE(Z3:Z51)"
$ws.Range("AA54").Formula = "=AVERAGE(AA3:AA51)"
$ws.Range("AB54").Formula = "=AVERAGE(AB3:AB51)"

# Row 56 flag column (mirrors the 0/1 "good init" flags used by the other
# packet blocks).
$ws.Range("Z56").Value = 0

# --- View state ----------------------------------------------------
# Scroll the visible window down to match where the new data was entered,
# and leave the selection on the cell below the new block.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G57").Select()
